$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the "time out" values for a few logged days (D column) - the
# dependent duration formulas in column E recalc automatically.
$ws.Range("D8").Value2 = 0.64583333333333337
$ws.Range("D13").Value2 = 0.61458333333333337
$ws.Range("D17").Value2 = 0.60416666666666663

# Update the active selection on the sheet.
$ws.Range("G5").Select()
